$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 721.5
$ws.Range("I29").Value = 610.2857
$ws.Range("J29").Value = 1500
$ws.Range("K29").Value = 1830.8571
$ws.Range("L29").Value = 4500
$ws.Range("M29").Value = -1549.8571
$ws.Range("N29").Value = -5062
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 898.04346
$ws.Range("I97").Value = 721.0714
$ws.Range("J97").Value = 1173.3334
$ws.Range("K97").Value = 721.0714
$ws.Range("L97").Value = 1173.3334
$ws.Range("M97").Value = -225.0714
$ws.Range("N97").Value = -2165.3334
$ws.Range("H102").Value = 1107.5
$ws.Range("I102").Value = 1107.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1107.5
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = 514.5
$ws.Range("H122").Value = 2388
$ws.Range("J122").Value = 2145.4
$ws.Range("L122").Value = 6436.200000000001
$ws.Range("N122").Value = -11336.2
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 920.0526
$ws.Range("I94").Value = 620.85
$ws.Range("J94").Value = 1252.5
$ws.Range("K94").Value = 620.85
$ws.Range("L94").Value = 1252.5
$ws.Range("M94").Value = -169.85
$ws.Range("N94").Value = -2154.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 35.041668
$ws.Range("J7").Value = 46.18182
$ws.Range("L7").Value = 46.18182
$ws.Range("N7").Value = -272.18182
$ws.Range("H86").Value = 2913.9375
$ws.Range("I86").Value = 2323.3333
$ws.Range("J86").Value = 4041.4546
$ws.Range("K86").Value = 2323.3333
$ws.Range("L86").Value = 4041.4546
$ws.Range("M86").Value = -1200.3333
$ws.Range("N86").Value = -6287.4546
$ws.Range("H89").Value = 2913.9375
$ws.Range("I89").Value = 2323.3333
$ws.Range("J89").Value = 4041.4546
$ws.Range("K89").Value = 11616.6665
$ws.Range("L89").Value = 20207.273
$ws.Range("M89").Value = -6000.666499999999
$ws.Range("N89").Value = -31439.273
$ws.Range("H105").Value = 7143.5
$ws.Range("I105").Value = 4572.2
$ws.Range("K105").Value = 4572.2
$ws.Range("M105").Value = -2825.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 223338.78
$ws.Range("I22").Value = 666750.3
$ws.Range("J22").Value = 1633
$ws.Range("K22").Value = 2000250.9
$ws.Range("L22").Value = 4899
$ws.Range("M22").Value = -2000081.9
$ws.Range("N22").Value = -5237
$ws.Range("H27").Value = 223338.78
$ws.Range("I27").Value = 666750.3
$ws.Range("J27").Value = 1633
$ws.Range("K27").Value = 2000250.9
$ws.Range("L27").Value = 4899
$ws.Range("M27").Value = -2000148.9
$ws.Range("N27").Value = -5103
$ws.Range("H34").Value = 954.8261
$ws.Range("J34").Value = 1161.1765
$ws.Range("L34").Value = 3483.5295
$ws.Range("N34").Value = -3651.5295
$ws.Range("I46").Value = 650
$ws.Range("J46").Value = 2150
$ws.Range("K46").Value = 1950
$ws.Range("L46").Value = 6450
$ws.Range("M46").Value = -1859
$ws.Range("N46").Value = -6632
$ws.Range("H61").Value = 80
$ws.Range("I61").Value = 80
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 240
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -25
$ws.Range("H88").Value = 4875
$ws.Range("J88").Value = 4875
$ws.Range("L88").Value = 14625
$ws.Range("N88").Value = -15481
$ws.Range("H91").Value = 4875
$ws.Range("J91").Value = 4875
$ws.Range("L91").Value = 14625
$ws.Range("N91").Value = -17589
$ws.Range("H130").Value = 2671.4285
$ws.Range("I130").Value = 700
$ws.Range("J130").Value = 3000
$ws.Range("K130").Value = 2100
$ws.Range("L130").Value = 9000
$ws.Range("M130").Value = 2920
$ws.Range("N130").Value = -19040
$ws.Range("H132").Value = 29860.5
$ws.Range("I132").Value = 834.1667
$ws.Range("J132").Value = 43257.27
$ws.Range("K132").Value = 7507.5003
$ws.Range("L132").Value = 389315.43
$ws.Range("M132").Value = -4977.5003
$ws.Range("N132").Value = -394375.43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 12009.25
$ws.Range("J26").Value = 12009.25
$ws.Range("L26").Value = 12009.25
$ws.Range("N26").Value = -12569.25
$ws.Range("H39").Value = 14800
$ws.Range("J39").Value = 14800
$ws.Range("L39").Value = 14800
$ws.Range("N39").Value = -15864
$ws.Range("H50").Value = 12009.25
$ws.Range("J50").Value = 12009.25
$ws.Range("L50").Value = 12009.25
$ws.Range("N50").Value = -13005.25
$ws.Range("H52").Value = 12000
$ws.Range("J52").Value = 12000
$ws.Range("L52").Value = 12000
$ws.Range("N52").Value = -12518
$ws.Range("H70").Value = 4001.6
$ws.Range("I70").Value = 4002.6667
$ws.Range("K70").Value = 4002.6667
$ws.Range("M70").Value = -3732.6667
$ws.Range("H73").Value = 4001.6
$ws.Range("I73").Value = 4002.6667
$ws.Range("K73").Value = 4002.6667
$ws.Range("M73").Value = -3066.6667
$ws.Range("H97").Value = 995.2
$ws.Range("I97").Value = 936
$ws.Range("J97").Value = 1133.3334
$ws.Range("K97").Value = 936
$ws.Range("L97").Value = 1133.3334
$ws.Range("M97").Value = -440
$ws.Range("N97").Value = -2125.3334
$ws.Range("H122").Value = 3819
$ws.Range("I122").Value = 1208.25
$ws.Range("J122").Value = 7300
$ws.Range("K122").Value = 3624.75
$ws.Range("L122").Value = 21900
$ws.Range("M122").Value = -1174.75
$ws.Range("N122").Value = -26800
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 4347.2856
$ws.Range("I61").Value = 4347.2856
$ws.Range("K61").Value = 4347.2856
$ws.Range("M61").Value = -4055.2856
$ws.Range("H96").Value = 1651
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1651
$ws.Range("K96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("M96").Value = 1651
$ws.Range("N96").Value = -4397
$ws.Range("H126").Value = 961.71875
$ws.Range("I126").Value = 865.1818
$ws.Range("J126").Value = 1174.1
$ws.Range("K126").Value = 2595.5454
$ws.Range("L126").Value = 3522.3
$ws.Range("M126").Value = -125.5454
$ws.Range("N126").Value = -8462.299999999999
